$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected "ParticipantsTab" query text (fix from Yizhen for \omic\ in cds corrected
# test cases in library selection cds): the WHERE clause is lower-cased, a "WITH p"
# line is inserted right after it (dropping samp/f/g/diag from scope before the
# re-MATCH on study/sample), "RETURN" loses its trailing space, and the final
# "LIMIT" is lower-cased.
$newQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
where g.library_selection in ['rRNA Depletion']
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id limit 100
'@

# Cell B2 is the only cell holding the old (incorrect) query text, so replacing its
# value in place corrects it; C2/C3/C4 keep pointing at the unrelated StatQuery
# ("CALL{ ... }") text, which is unaffected.
$ws.Range("B2").Value = $newQuery

# Update the active selection from B4 to C3, as recorded in sheetView/selection.
$ws.Range("C3").Select()
